$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "Region" header in column O, matching the style of the other headers (column N)
$ws.Range("N1").Copy()
$ws.Range("O1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("O1").Value2 = "Region"

# Lookup table mapping State abbreviation/name (column N) to Region
$regionMap = @{
    "HI" = "Other"
    "Unknown" = "Other"
    "Guam" = "Other"
    "American Samoa" = "Other"
    "Marshall Islands" = "Other"
    "Pacific Ocean" = "Other"
    "Bermuda" = "Other"
    "ME" = "East Coast"
    "NH" = "East Coast"
    "MA" = "East Coast"
    "RI" = "East Coast"
    "CT" = "East Coast"
    "NY" = "East Coast"
    "NJ" = "East Coast"
    "PA" = "Other"
    "DE" = "East Coast"
    "MD" = "East Coast"
    "DC" = "Other"
    "VA" = "East Coast"
    "NC" = "East Coast"
    "SC" = "East Coast"
    "GA" = "East Coast"
    "FL" = "Gulf Coast"
    "AL" = "Gulf Coast"
    "MS" = "Gulf Coast"
    "LA" = "Gulf Coast"
    "TX" = "Gulf Coast"
    "CA" = "West Coast"
    "OR" = "West Coast"
    "WA" = "West Coast"
    "AK" = "Other"
    "VI" = "Other"
    "PR" = "Other"
}

# Populate column O (Region) for each data row based on column N (State)
for ($row = 2; $row -le 146; $row++) {
    $state = $ws.Cells.Item($row, 14).Value2
    $region = $regionMap[$state]
    $ws.Cells.Item($row, 15).Value2 = $region
}
